$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row -> (F value, G value) to update, per the diff.
$updates = @(
    @{ Row = 12;  F = "0,1,2,3"; G = "aoe" },
    @{ Row = 23;  F = "0,1,2,3"; G = "random3" },
    @{ Row = 26;  F = "0,1,2,3"; G = "aoe,ally" },
    @{ Row = 27;  F = "0,1,2,3"; G = "random" },
    @{ Row = 44;  F = "0,1,2,3"; G = "aoe" },
    @{ Row = 45;  F = "0,1,2,3"; G = "aoe" },
    @{ Row = 49;  F = "0,1,2,3"; G = "aoe" },
    @{ Row = 53;  F = "0,1,2,3"; G = "aoe" },
    @{ Row = 62;  F = "0,1,2,3"; G = "aoe" },
    @{ Row = 67;  F = "0,1,2,3"; G = "random2" },
    @{ Row = 70;  F = "0,1,2,3"; G = "random3" },
    @{ Row = 80;  F = "0,1,2,3"; G = "aoe" },
    @{ Row = 83;  F = "0,1,2,3"; G = "random2" },
    @{ Row = 93;  F = "0,1,2,3"; G = "aoe" },
    @{ Row = 101; F = "0,1,2,3"; G = "aoe" },
    @{ Row = 140; F = "0,1,2,3"; G = "aoe,ally" },
    @{ Row = 143; F = "0,1,2,3"; G = "aoe" },
    @{ Row = 144; F = "0,1,2,3"; G = "aoe" },
    @{ Row = 159; F = "0,1,2,3"; G = "aoe" },
    @{ Row = 175; F = "0,1,2,3"; G = "aoe" },
    @{ Row = 179; F = "0,1,2,3"; G = "aoe,ally" },
    @{ Row = 211; F = "0,1,2,3"; G = "aoe" },
    @{ Row = 214; F = "0,1,2,3"; G = "aoe" },
    @{ Row = 231; F = "0,1,2,3"; G = "aoe" },
    @{ Row = 238; F = "0,1,2,3"; G = "aoe,ally" },
    @{ Row = 245; F = "0,1,2,3"; G = "aoe" },
    @{ Row = 246; F = "0,1,2,3"; G = "aoe" },
    @{ Row = 247; F = "0,1,2,3"; G = "aoe" },
    @{ Row = 257; F = "0,1,2,3"; G = "aoe,ally" },
    @{ Row = 258; F = "0,1,2,3"; G = "aoe" },
    @{ Row = 260; F = "0,1,2,3"; G = "aoe" }
)

foreach ($u in $updates) {
    $r = $u.Row
    $ws.Range("F$r").Value = $u.F
    $ws.Range("G$r").Value = $u.G
}
